# Update F-column values (想去人数 / interest counts) across all 4 sheets
# as per the commit 'Update gh-pages to output generated at 456a3b4'
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1886
$ws.Range("F4").Value = 72
$ws.Range("F5").Value = 744
$ws.Range("F6").Value = 116
$ws.Range("F7").Value = 501
$ws.Range("F8").Value = 904
$ws.Range("F9").Value = 1607
$ws.Range("F10").Value = 1274
$ws.Range("F11").Value = 1539
$ws.Range("F12").Value = 66
$ws.Range("F13").Value = 1524
$ws.Range("F14").Value = 343
$ws.Range("F15").Value = 1690
$ws.Range("F16").Value = 807
$ws.Range("F17").Value = 1112
$ws.Range("F18").Value = 371
$ws.Range("F19").Value = 57
$ws.Range("F21").Value = 1743
$ws.Range("F22").Value = 231
$ws.Range("F23").Value = 822
$ws.Range("F24").Value = 3
$ws.Range("F25").Value = 560
$ws.Range("F26").Value = 1230
$ws.Range("F27").Value = 325365
$ws.Range("F28").Value = 1061
$ws.Range("F29").Value = 78
$ws.Range("F30").Value = 575
$ws.Range("F31").Value = 1340
$ws.Range("F32").Value = 1156
$ws.Range("F34").Value = 8
$ws.Range("F35").Value = 1155
$ws.Range("F36").Value = 1103
$ws.Range("F37").Value = 275
$ws.Range("F38").Value = 79
$ws.Range("F39").Value = 883
$ws.Range("F40").Value = 1688
$ws.Range("F42").Value = 120
$ws.Range("F44").Value = 2042
$ws.Range("F45").Value = 91
$ws.Range("F46").Value = 833
$ws.Range("F47").Value = 804
$ws.Range("F49").Value = 30

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 49
$ws.Range("F3").Value = 50
$ws.Range("F9").Value = 2594
$ws.Range("F10").Value = 1221
$ws.Range("F11").Value = 410
$ws.Range("F12").Value = 727
$ws.Range("F13").Value = 257
$ws.Range("F14").Value = 37
$ws.Range("F18").Value = 465
$ws.Range("F21").Value = 318
$ws.Range("F22").Value = 0
$ws.Range("F24").Value = 4
$ws.Range("F27").Value = 192
$ws.Range("F28").Value = 250
$ws.Range("F30").Value = 217
$ws.Range("F33").Value = 58
$ws.Range("F34").Value = 10
$ws.Range("F36").Value = 11
$ws.Range("F37").Value = 2
$ws.Range("F38").Value = 184
$ws.Range("F42").Value = 62
$ws.Range("F43").Value = 62
$ws.Range("F45").Value = 140
$ws.Range("F46").Value = 66

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 120
$ws.Range("F4").Value = 269
$ws.Range("F5").Value = 2896
$ws.Range("F6").Value = 4670
$ws.Range("F9").Value = 584
$ws.Range("F10").Value = 752
$ws.Range("F11").Value = 476
$ws.Range("F12").Value = 383
$ws.Range("F13").Value = 1126
$ws.Range("F14").Value = 298
$ws.Range("F15").Value = 710

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1886
$ws.Range("F3").Value = 269
$ws.Range("F5").Value = 4670
$ws.Range("F6").Value = 752
$ws.Range("F7").Value = 476
$ws.Range("F8").Value = 383
$ws.Range("F9").Value = 383
$ws.Range("F10").Value = 1126
$ws.Range("F11").Value = 501
$ws.Range("F12").Value = 904
$ws.Range("F13").Value = 1221
$ws.Range("F14").Value = 1607
$ws.Range("F15").Value = 1274
$ws.Range("F16").Value = 1539
$ws.Range("F17").Value = 1524
$ws.Range("F18").Value = 257
$ws.Range("F20").Value = 1690
$ws.Range("F21").Value = 1112
$ws.Range("F22").Value = 371
$ws.Range("F23").Value = 710
$ws.Range("F24").Value = 710
$ws.Range("F25").Value = 1743
$ws.Range("F26").Value = 231
$ws.Range("F27").Value = 822
$ws.Range("F28").Value = 560
$ws.Range("F29").Value = 1230
$ws.Range("F30").Value = 318
$ws.Range("F31").Value = 1061
$ws.Range("F32").Value = 78
$ws.Range("F33").Value = 1156
$ws.Range("F35").Value = 8
$ws.Range("F36").Value = 1155
$ws.Range("F37").Value = 4
$ws.Range("F38").Value = 1103
$ws.Range("F39").Value = 275
$ws.Range("F40").Value = 883
$ws.Range("F42").Value = 1688
$ws.Range("F44").Value = 120
$ws.Range("F45").Value = 2042
$ws.Range("F46").Value = 91
$ws.Range("F47").Value = 833
$ws.Range("F48").Value = 62
$ws.Range("F49").Value = 804
$ws.Range("F52").Value = 66
